$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 195-196, pushing the old rows 195..214 down to 197..216.
$ws.Range("A195:A196").EntireRow.Insert()

# New row 195: Pera / Packham's Triumph / Especial, fecha 2022-08-10 (serial 44783)
$ws.Range("A195").Value = 7
$ws.Range("B195").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C195").Value = "Ñuble"
$ws.Range("D195").Value = 44783
$ws.Range("E195").Value = 16
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100104
$ws.Range("H195").Value = "Frutos de pepita"
$ws.Range("I195").Value = 100104005
$ws.Range("J195").Value = "Pera"
$ws.Range("K195").Value = "Packham's Triumph"
$ws.Range("L195").Value = "Especial"
$ws.Range("M195").Value = 40
$ws.Range("N195").Value = 10000
$ws.Range("O195").Value = 10000
$ws.Range("P195").Value = 10000
$ws.Range("Q195").Value = "$/caja 16 kilos empedrada"
$ws.Range("R195").Value = "Provincia de Curicó"
$ws.Range("S195").Value = 625
$ws.Range("T195").Value = 16

# New row 196: Pera / Packham's Triumph / Primera, fecha 2022-08-10 (serial 44783)
$ws.Range("A196").Value = 7
$ws.Range("B196").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C196").Value = "Ñuble"
$ws.Range("D196").Value = 44783
$ws.Range("E196").Value = 16
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100104
$ws.Range("H196").Value = "Frutos de pepita"
$ws.Range("I196").Value = 100104005
$ws.Range("J196").Value = "Pera"
$ws.Range("K196").Value = "Packham's Triumph"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 80
$ws.Range("N196").Value = 6500
$ws.Range("O196").Value = 9000
$ws.Range("P196").Value = 7750
$ws.Range("Q196").Value = "$/caja 16 kilos empedrada"
$ws.Range("R196").Value = "Provincia de Curicó"
$ws.Range("S196").Value = 484
$ws.Range("T196").Value = 16
